$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 changes
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("BD2").Value = 151

# Row 4 changes
$ws.Range("G4").Value = 2.45
$ws.Range("I4").Value = 2.75
$ws.Range("J4").Value = 3.1
$ws.Range("L4").Value = 3.4
$ws.Range("N4").Value = 10
$ws.Range("O4").Value = 1.3
$ws.Range("P4").Value = 3.4
$ws.Range("Q4").Value = 1.98
$ws.Range("R4").Value = 1.83
$ws.Range("S4").Value = 1.4
$ws.Range("T4").Value = 2.75
$ws.Range("U4").Value = 1.73
$ws.Range("V4").Value = 2
$ws.Range("X4").Value = 12
$ws.Range("Z4").Value = 23
$ws.Range("AB4").Value = 29
$ws.Range("AC4").Value = 10
$ws.Range("AG4").Value = 9
$ws.Range("AH4").Value = 13
$ws.Range("AI4").Value = 11
$ws.Range("AJ4").Value = 29
$ws.Range("AK4").Value = 23
$ws.Range("AM4").Value = 201
$ws.Range("AN4").Value = 4.5
$ws.Range("AO4").Value = 13
$ws.Range("AQ4").Value = 41
$ws.Range("AT4").Value = 2.75
$ws.Range("AW4").Value = 4.75
